# Update lattice-multiplication practice grid: 15 cells (5 rows x 3 cols)
# each cell's 5-line content (problem, factors, separator, two partial digits)
# is replaced in place via the Tables/Cell Range API so w:br line breaks and
# run formatting (sz=32) are regenerated faithfully.
$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$cell = $tbl.Cell(1, 1)
$cell.Range.Text = "48 x 94`v  9    4`v  ----`v4|    |`v8|    |"
$cell = $tbl.Cell(1, 2)
$cell.Range.Text = "25 x 65`v  6    5`v  ----`v2|    |`v5|    |"
$cell = $tbl.Cell(1, 3)
$cell.Range.Text = "57 x 25`v  2    5`v  ----`v5|    |`v7|    |"
$cell = $tbl.Cell(2, 1)
$cell.Range.Text = "21 x 84`v  8    4`v  ----`v2|    |`v1|    |"
$cell = $tbl.Cell(2, 2)
$cell.Range.Text = "83 x 59`v  5    9`v  ----`v8|    |`v3|    |"
$cell = $tbl.Cell(2, 3)
$cell.Range.Text = "93 x 56`v  5    6`v  ----`v9|    |`v3|    |"
$cell = $tbl.Cell(3, 1)
$cell.Range.Text = "86 x 56`v  5    6`v  ----`v8|    |`v6|    |"
$cell = $tbl.Cell(3, 2)
$cell.Range.Text = "63 x 54`v  5    4`v  ----`v6|    |`v3|    |"
$cell = $tbl.Cell(3, 3)
$cell.Range.Text = "64 x 39`v  3    9`v  ----`v6|    |`v4|    |"
$cell = $tbl.Cell(4, 1)
$cell.Range.Text = "47 x 54`v  5    4`v  ----`v4|    |`v7|    |"
$cell = $tbl.Cell(4, 2)
$cell.Range.Text = "95 x 19`v  1    9`v  ----`v9|    |`v5|    |"
$cell = $tbl.Cell(4, 3)
$cell.Range.Text = "19 x 56`v  5    6`v  ----`v1|    |`v9|    |"
$cell = $tbl.Cell(5, 1)
$cell.Range.Text = "86 x 35`v  3    5`v  ----`v8|    |`v6|    |"
$cell = $tbl.Cell(5, 2)
$cell.Range.Text = "91 x 11`v  1    1`v  ----`v9|    |`v1|    |"
$cell = $tbl.Cell(5, 3)
$cell.Range.Text = "35 x 59`v  5    9`v  ----`v3|    |`v5|    |"

Write-Output "Lattice multiplication grid updated."
